# =====================================================================
# Adds a new "2022-Q1" sheet (fund holding detail) before the "总计"
# (summary) sheet, and updates "总计" with a new leading row for 2022-Q1.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Capture the existing "总计" sheet's data (header text) before we
#    recreate it, and remember a style-2 source cell to copy formats
#    from (column-A index cells / header cells all use the same bold+
#    bordered+centered style).
# ---------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item(1).Cells.Item(2, 1)   # a cell already using style "2"

$oldTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Delete()

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q1" sheet right after "2021-Q4", then the
#    "总计" sheet right after that — so sheet order/ids come out as
#    ..., 2021-Q4, 2022-Q1, 总计.
# ---------------------------------------------------------------------
$q4sheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newQ1 = $wb.Worksheets.Add($null, $q4sheet)
$newQ1.Name = "2022-Q1"

$newTotal = $wb.Worksheets.Add($null, $newQ1)
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# Helper: write a cell as genuine text (keeps leading zeros / decimal
# strings verbatim instead of Excel's automatic number coercion), then
# drop back to the default "Normal" style so no stray number format
# lingers on the cell.
# ---------------------------------------------------------------------
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-IndexCell($cell, $num) {
    $cell.Value = $num
    $styleSrc.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------
# 3. Populate "2022-Q1" — header row + 46 fund rows.
# ---------------------------------------------------------------------
$ws1 = $newQ1

Set-TextCell $ws1.Cells.Item(1, 2) "基金代码"
Set-TextCell $ws1.Cells.Item(1, 3) "基金名称"
Set-TextCell $ws1.Cells.Item(1, 4) "基金规模"
Set-TextCell $ws1.Cells.Item(1, 5) "股票总仓位"
Set-TextCell $ws1.Cells.Item(1, 6) "仓位占比"
Set-TextCell $ws1.Cells.Item(1, 7) "持有市值(亿元)"
Set-TextCell $ws1.Cells.Item(1, 8) "仓位排名"
$styleSrc.Copy()
$ws1.Range("B1:H1").PasteSpecial(-4122)

$fundData = @(
    ,("0","005267","嘉实价值精选股票","65.04","92.05","7.82","5.0861","4")
    ,("1","006567","中泰星元价值优选灵活配置混合","44.13","81.95","9.31","4.1085","2")
    ,("2","010273","嘉实价值长青混合A","47.29","86.91","6.49","3.0691","5")
    ,("3","012533","嘉实价值驱动一年持有期混合型证券投资基金A","50.84","91.01","5.60","2.8470","5")
    ,("4","011518","嘉实价值臻选混合型证券投资基金","35.83","89.89","7.09","2.5403","3")
    ,("5","012344","嘉实领先优势混合型证券投资基金A","68.76","82.99","3.02","2.0766","8")
    ,("6","007549","中泰开阳价值优选灵活配置混合A","20.45","91.42","9.57","1.9571","4")
    ,("7","013776","中泰兴为价值精选混合A","20.31","85.34","9.55","1.9396","2")
    ,("8","070019","嘉实价值优势混合","27.51","93.14","6.95","1.9119","4")
    ,("9","010190","嘉实价值发现三个月定期开放混合","33.17","93.94","5.50","1.8244","5")
    ,("10","006624","中泰玉衡价值优选混合","17.75","81.95","9.27","1.6454","2")
    ,("11","001878","嘉实沪港深精选股票","23.17","93.29","5.41","1.2535","3")
    ,("12","012001","中泰星宇价值成长一年封闭运作混合型证券投资基金A","12.14","91.96","9.71","1.1788","4")
    ,("13","013777","中泰兴为价值精选混合C","8.71","85.34","9.55","0.8318","2")
    ,("14","070003","嘉实稳健混合","24.15","70.47","3.27","0.7897","9")
    ,("15","001044","嘉实新消费股票","8.92","80.25","8.23","0.7341","2")
    ,("16","004355","嘉实丰和灵活配置混合","9.22","85.59","6.71","0.6187","3")
    ,("17","010728","中泰兴诚价值一年持有期混合A","9.94","88.08","5.09","0.5059","6")
    ,("18","011437","中泰开阳价值优选灵活配置混合C","5.02","91.42","9.57","0.4804","4")
    ,("19","671010","西部利得策略优选混合A","6.79","93.07","5.06","0.3436","10")
    ,("20","010274","嘉实价值长青混合C","3.97","86.91","6.49","0.2577","5")
    ,("21","011521","鹏扬景源一年持有期混合A","33.78","21.39","0.74","0.2500","6")
    ,("22","920002","中金精选股票A","3.40","86.43","6.96","0.2366","6")
    ,("23","012002","中泰星宇价值成长一年封闭运作混合型证券投资基金C","2.13","91.96","9.71","0.2068","4")
    ,("24","012534","嘉实价值驱动一年持有期混合型证券投资基金C","3.20","91.01","5.60","0.1792","5")
    ,("25","000574","宝盈新价值灵活配置混合A","4.32","84.13","4.00","0.1728","10")
    ,("26","003715","宝盈消费主题灵活配置混合","2.74","83.48","4.01","0.1099","9")
    ,("27","001577","嘉实低价策略股票","1.53","85.79","6.18","0.0946","3")
    ,("28","010729","中泰兴诚价值一年持有期混合C","1.72","88.08","5.09","0.0875","6")
    ,("29","005265","博时厚泽回报灵活配置混合A","2.19","91.85","3.34","0.0731","6")
    ,("30","001707","诺安高端制造股票","1.35","91.74","5.05","0.0682","4")
    ,("31","011060","西部利得策略优选混合C","1.32","93.07","5.06","0.0668","10")
    ,("32","000963","兴业多策略灵活配置混合","2.07","75.34","2.92","0.0604","8")
    ,("33","012153","博时研究慧选混合型证券投资基金A","1.63","75.28","3.06","0.0499","9")
    ,("34","011522","鹏扬景源一年持有期混合C","4.58","21.39","0.74","0.0339","6")
    ,("35","005041","人保研究精选混合A","1.33","81.87","2.47","0.0329","2")
    ,("36","009766","安信平稳双利3个月持有期混合A","2.33","39.45","1.17","0.0273","10")
    ,("37","012345","嘉实领先优势混合型证券投资基金C","0.82","82.99","3.02","0.0248","8")
    ,("38","005266","博时厚泽回报灵活配置混合C","0.64","91.85","3.34","0.0214","6")
    ,("39","920922","中金精选股票C","0.14","86.43","6.96","0.0097","6")
    ,("40","012154","博时研究慧选混合型证券投资基金C","0.21","75.28","3.06","0.0064","9")
    ,("41","007574","宝盈新价值灵活配置混合C","0.09","84.13","4.00","0.0036","10")
    ,("42","009767","安信平稳双利3个月持有期混合C","0.26","39.45","1.17","0.0030","10")
    ,("43","750005","安信平稳增长混合A","0.08","65.16","2.72","0.0022","9")
    ,("44","005042","人保研究精选混合C","0.03","81.87","2.47","0.0007","2")
    ,("45","002035","安信平稳增长混合C","0.00","65.16","2.72","0","9")
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $row = $i + 2
    $rec = $fundData[$i]
    $idxVal = $rec[0]
    Set-IndexCell $ws1.Cells.Item($row, 1) $idxVal
    Set-TextCell $ws1.Cells.Item($row, 2) $rec[1]
    Set-TextCell $ws1.Cells.Item($row, 3) $rec[2]
    Set-TextCell $ws1.Cells.Item($row, 4) $rec[3]
    Set-TextCell $ws1.Cells.Item($row, 5) $rec[4]
    Set-TextCell $ws1.Cells.Item($row, 6) $rec[5]
    if ($row -eq 47) {
        $ws1.Cells.Item($row, 7).Value = 0
    } else {
        Set-TextCell $ws1.Cells.Item($row, 7) $rec[6]
    }
    $rankVal = $rec[7]
    $ws1.Cells.Item($row, 8).Value = $rankVal
}

# ---------------------------------------------------------------------
# 4. Populate "总计" — header row + 6 summary rows (2022-Q1 first, then
#    the previously-existing quarters shifted down by one).
# ---------------------------------------------------------------------
$ws2 = $newTotal

Set-TextCell $ws2.Cells.Item(1, 2) "日期"
Set-TextCell $ws2.Cells.Item(1, 3) "持有数量(只)"
Set-TextCell $ws2.Cells.Item(1, 4) "持有市值(亿元)"
$styleSrc.Copy()
$ws2.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    ,("2022-Q1", 46, 37.82)
    ,("2021-Q4", 45, 25.44)
    ,("2021-Q3", 40, 26.53)
    ,("2021-Q2", 34, 21.91)
    ,("2021-Q1", 96, 47.57)
    ,("2020-Q4", 52, 23.3)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $row = $i + 2
    $rec = $totalRows[$i]
    Set-IndexCell $ws2.Cells.Item($row, 1) $i
    Set-TextCell $ws2.Cells.Item($row, 2) $rec[0]
    $ws2.Cells.Item($row, 3).Value = $rec[1]
    $ws2.Cells.Item($row, 4).Value = $rec[2]
}
